# Update template metadata for new database
# Applies the edits described in the commit: restructures the
# "SwateTemplateMetadata" sheet (drop Docslink row, expand Tags row to
# 3 values, rename "Authors Roles" -> "Authors Role" and add an
# "Authors ORCID" field), and adds the reviewer's threaded comments that
# document each metadata field.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("SwateTemplateMetadata")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Remove the obsolete "Docslink" row (old row 5). Everything below
#    shifts up by one row.
# ---------------------------------------------------------------------
$ws.Rows.Item(5).Delete()

# ---------------------------------------------------------------------
# 2) Tags block now carries 3 example values instead of 1 (row 12),
#    and the accession/source-ref rows (13/14) grow B:D to match.
# ---------------------------------------------------------------------

# Header row "#TAGS list" (row 11): extend the banner fill across C:D.
$ws.Range("B11").Copy() | Out-Null
$ws.Range("C11:D11").PasteSpecial($xlPasteFormats) | Out-Null

# "Tags" row (row 12): add the two extra example tags + matching style,
# and bump the row height to fit the wrapped header text.
$ws.Range("B12").Copy() | Out-Null
$ws.Range("C12:D12").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C12").Value = "growth factors"
$ws.Range("D12").Value = "conditions"
$ws.Rows.Item(12).RowHeight = 30

# "Tags Term Accession Number" / "Tags Term Source REF" rows (13/14).
$ws.Range("B13").Copy() | Out-Null
$ws.Range("C13:D13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B14").Copy() | Out-Null
$ws.Range("C14:D14").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Authors Roles -> Authors Role, plus a new "Authors ORCID" field.
#    After the Docslink deletion the layout is:
#      23 Authors Affiliation
#      24 #AUTHORS ROLES list   (section header)
#      25 Authors Roles
#      26 Authors Roles Term Accession Number
#      27 Authors Roles Term Source REF   (last row of the sheet)
# ---------------------------------------------------------------------

# Drop the "#AUTHORS ROLES list" section header entirely.
$ws.Rows.Item(24).Delete()

# Make room for the new "Authors ORCID" row right above "Authors Roles".
$ws.Rows.Item(24).Insert()

# New row 24: Authors ORCID, styled like the row above it (Affiliation).
$ws.Range("A23:B23").Copy() | Out-Null
$ws.Range("A24:B24").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A24").Value = "Authors ORCID"
$ws.Range("B24").Value = ""

# "Authors Roles" (now row 25) loses its "first row under a header"
# shading (style 8 -> 6) now that the header above it is gone.
$ws.Range("B26").Copy() | Out-Null
$ws.Range("B25").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A25").Value = "Authors Role"

# Rename the two remaining rows.
$ws.Range("A26").Value = "Authors Role Term Accession Number"
$ws.Range("A27").Value = "Authors Role Term Source REF"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Column / view tweaks that came along with the restructuring.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 15.140625
$ws.Range("C7").Select()

# ---------------------------------------------------------------------
# 5) Threaded comments documenting each metadata field (added by the
#    new reviewer).
# ---------------------------------------------------------------------
$ws.Range("A2").AddCommentThreaded("The name of the Swate template.") | Out-Null
$ws.Range("A3").AddCommentThreaded("The current version of this template in SemVer notation.") | Out-Null
$ws.Range("A4").AddCommentThreaded("The description of this template. Use few sentences for succinctness.") | Out-Null
$ws.Range("A5").AddCommentThreaded('The name of the template associated organisation. "DataPLANT" will trigger the "DataPLANT" batch of honor for the template.') | Out-Null
$ws.Range("A6").AddCommentThreaded("The name of the Swate annotation table in the workbook of the template's excel file.") | Out-Null
$ws.Range("A7").AddCommentThreaded("A list of all ERs (endpoint repositories) targeted with this template. ERs are realized as Terms: <term ref here>") | Out-Null
$ws.Range("A11").AddCommentThreaded("A list of all tags associated with this template. Tags are realized as Terms: <term ref here>") | Out-Null
$ws.Range("A15").AddCommentThreaded("The author(s) of this template.") | Out-Null
